$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @(3.286832544864788, 1.655778082260271, 3.537761648806719, 0.4942365360607697, 1, 8.974608811992548)
    3 = @(0.6606524410359556, 0.306821227259698, 0.1494219747398047, 0.4942365360607697, 1, 1.611132179096228)
    4 = @(0.2917716402565462, 0.04071648406533734, 3.537761648806719, 10.19245300693656, 0, 14.06270278006516)
    5 = @(3.286832544864788, 1.655778082260271, 0.7527432677738641, 0.4942365360607697, 1, 6.189590430959694)
    6 = @(3.286832544864788, 1.655778082260271, 3.537761648806719, 10.19245300693656, 0, 18.67282528286833)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]
    $ws.Cells.Item($row, 3).Value = $vals[1]
    $ws.Cells.Item($row, 4).Value = $vals[2]
    $ws.Cells.Item($row, 5).Value = $vals[3]
    $ws.Cells.Item($row, 6).Value = $vals[4]
    $ws.Cells.Item($row, 7).Value = $vals[5]
}
